$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old dummy row 3 content
$ws.Range("A3").ClearContents()

# Update header row
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "pi_type"
$ws.Range("C1").Value = "cs_pi_name"
$ws.Range("D1").Value = "other_pi_name"

# Update data row
$ws.Range("A2").Value = "Grant for testing"
$ws.Range("B2").Value = "OTHER_PI"
$ws.Range("C2").Value = "test"
$ws.Range("D2").Value = "othertest"

$ws.Range("C2").Select()
